$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores values as text in the source data (e.g. big
# numbers use "." as a thousands separator, like "63.194.27", which is not a
# valid number anyway). Whenever a replacement price happens to look like a
# normal decimal number (e.g. "546.21"), mark the cell as Text first so Excel
# keeps storing it as a string instead of silently converting it to a number.

$ws.Range("D2").Value = "63.194.27"
$ws.Range("E2").Value = "  -3.11%  "

$ws.Range("D3").Value = "3.078.04"
$ws.Range("E3").Value = "  -1.58%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.21"
$ws.Range("E5").Value = "  -2.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.50"
$ws.Range("E6").Value = "  -7.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "3.071.05"
$ws.Range("E8").Value = "  -1.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -2.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("E10").Value = "  -0.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.53"
$ws.Range("E11").Value = "  -3.58%  "

$ws.Range("E12").Value = "  -1.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "34.71"
$ws.Range("E13").Value = "  -5.35%  "

$ws.Range("E14").Value = "  -2.59%  "

$ws.Range("D15").Value = "3.585.41"
$ws.Range("E15").Value = "  -1.20%  "

$ws.Range("D16").Value = "63.293.06"
$ws.Range("E16").Value = "  -3.09%  "

$ws.Range("E17").Value = "  -0.93%  "

$ws.Range("D18").Value = "3.090.59"
$ws.Range("E18").Value = "  -1.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "500.21"
$ws.Range("E19").Value = "  -4.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.62"
$ws.Range("E20").Value = "  -2.26%  "

$ws.Range("E21").Value = "  -4.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.699"
$ws.Range("E22").Value = "  -1.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.15"
$ws.Range("E23").Value = "  -3.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.44"
$ws.Range("E24").Value = "  -2.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.26"
$ws.Range("E25").Value = "  -4.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.72"
$ws.Range("E27").Value = "  -2.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.31"
$ws.Range("E28").Value = "  -5.11%  "

$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("E30").Value = "  -9.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.28"
$ws.Range("E31").Value = "  -0.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  +1.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.51"
$ws.Range("E33").Value = "  -6.82%  "

$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "520.85"
$ws.Range("E34").Value = "  -9.70%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.21"
$ws.Range("E35").Value = "  +7.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.94"
$ws.Range("E36").Value = "  -2.80%  "

$ws.Range("E37").Value = "  -7.23%  "

$ws.Range("E38").Value = "  -7.46%  "

$ws.Range("D39").Value = "3.067.66"
$ws.Range("E39").Value = "  +0.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0789"
$ws.Range("E40").Value = "  -4.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.120"
$ws.Range("E41").Value = "  -2.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.67"
$ws.Range("E42").Value = "  -7.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.07"
$ws.Range("E43").Value = "  -2.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.255"
$ws.Range("E44").Value = "  -2.05%  "

$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.03"
$ws.Range("E46").Value = "  -7.47%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.61"
$ws.Range("E47").Value = "  +3.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.53"
$ws.Range("E48").Value = "  +63.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.23"
$ws.Range("E49").Value = "  -4.20%  "

$ws.Range("E50").Value = "  -1.96%  "

$ws.Range("D51").Value = "0.0₃0501"
$ws.Range("E51").Value = "  -6.37%  "

